# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) slide master -> "Integral" / "Red Violet"
#   ppt/theme/theme2.xml  -> bound to the notes master         -> "Office Theme" / "Office"
# The authored edit swaps the content of those two parts, so the slides end up
# themed with the default "Office Theme" palette while the notes master keeps
# the old "Integral" (Red Violet) palette.
#
# The PowerPoint object model reaches the deck's single editable theme through
# Design.SlideMaster.Theme (Master/NotesMaster/HandoutMaster all resolve to the
# same theme here), so we recolor that theme's 12 scheme slots to the "Office
# Theme" palette - this is the part of the swap that is externally visible
# (it is what slides actually render with).

function Hex2Rgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme palette, in clrScheme order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = Hex2Rgb($officeColors[$i])
}
